$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row (Fecha, Volumen, Precio minimo, Precio maximo,
# Precio promedio ponderado, Unidad de comercializacion, Origen,
# Precio $/Kg, Kg o Unidades)
$data = @{
    2  = @(44223, 80,  3500, 3800, 3688, "`$/paquete 2 kilos",      "Provincia de Diguillín", 1844, 2)
    3  = @(44208, 85,  3700, 4000, 3824, "`$/paquete 2 kilos",      "Provincia de Diguillín", 1912, 2)
    4  = @(44161, 50,  2800, 3000, 2900, "`$/paquete 2 kilos",      "Provincia de Diguillín", 1450, 2)
    5  = @(44704, 100, 6000, 6500, 6250, "`$/paquete 36 unidades",  "Región Metropolitana",   174,  36)
    6  = @(44225, 80,  3400, 3700, 3550, "`$/paquete 2 kilos",      "Provincia de Diguillín", 1775, 2)
    7  = @(44701, 120, 7000, 7500, 7250, "`$/paquete 36 unidades",  "Región Metropolitana",   201,  36)
    8  = @(44160, 43,  3500, 4000, 3709, "`$/paquete 36 unidades",  "Región Metropolitana",   103,  36)
    9  = @(44210, 105, 3500, 4000, 3714, "`$/paquete 2 kilos",      "Provincia de Diguillín", 1857, 2)
    10 = @(44215, 140, 3500, 4000, 3768, "`$/paquete 2 kilos",      "Provincia de Diguillín", 1884, 2)
    11 = @(44662, 200, 8000, 8500, 8250, "`$/paquete 36 unidades",  "Región Metropolitana",   229,  36)
    12 = @(44664, 200, 8000, 8500, 8250, "`$/paquete 36 unidades",  "Región Metropolitana",   229,  36)
    13 = @(44166, 70,  3500, 4000, 3679, "`$/paquete 36 unidades",  "Región Metropolitana",   102,  36)
    14 = @(44209, 150, 3500, 4000, 3767, "`$/paquete 2 kilos",      "Provincia de Diguillín", 1884, 2)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 4).Value  = $vals[0]   # D Fecha
    $ws.Cells.Item($row, 10).Value = $vals[1]   # J Volumen
    $ws.Cells.Item($row, 11).Value = $vals[2]   # K Precio minimo
    $ws.Cells.Item($row, 12).Value = $vals[3]   # L Precio maximo
    $ws.Cells.Item($row, 13).Value = $vals[4]   # M Precio promedio ponderado
    $ws.Cells.Item($row, 14).Value = $vals[5]   # N Unidad de comercializacion
    $ws.Cells.Item($row, 15).Value = $vals[6]   # O Origen
    $ws.Cells.Item($row, 16).Value = $vals[7]   # P Precio $/Kg
    $ws.Cells.Item($row, 17).Value = $vals[8]   # Q Kg o Unidades
}
